# "implement wifi raw rle data transfer"
# The WiFi timer (column B) prescaler is bumped from 20 to 1000; this
# cascades through the dependent formulas (B5, B7, B12) automatically on
# recalculation. Finally, leave the active selection on the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 1000

$ws.Range("B3").Select()
